$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 132
$ws.Cells.Item(2, 4).Value = 57947.84

# Row 3
$ws.Cells.Item(3, 3).Value = 53
$ws.Cells.Item(3, 4).Value = 48185

# Row 4
$ws.Cells.Item(4, 1).Value = 'SUCRIVOIRE'
$ws.Cells.Item(4, 3).Value = 54
$ws.Cells.Item(4, 4).Value = 45155
$ws.Cells.Item(4, 5).Value = 990

# Row 5
$ws.Cells.Item(5, 1).Value = 'BRVM - AUTRES SECTEURS'
$ws.Cells.Item(5, 3).Value = 66
$ws.Cells.Item(5, 4).Value = 44906.31
$ws.Cells.Item(5, 5).Value = 667.6

# Row 6
$ws.Cells.Item(6, 3).Value = 61
$ws.Cells.Item(6, 4).Value = 44080

# Row 8
$ws.Cells.Item(8, 3).Value = 66
$ws.Cells.Item(8, 4).Value = 39095

# Row 9
$ws.Cells.Item(9, 3).Value = 66
$ws.Cells.Item(9, 4).Value = 38030

# Row 10
$ws.Cells.Item(10, 3).Value = 66
$ws.Cells.Item(10, 4).Value = 32760

# Row 11
$ws.Cells.Item(11, 3).Value = 66
$ws.Cells.Item(11, 4).Value = 26445

# Row 12
$ws.Cells.Item(12, 3).Value = 66
$ws.Cells.Item(12, 4).Value = 24302.99

# Row 13
$ws.Cells.Item(13, 3).Value = 66
$ws.Cells.Item(13, 4).Value = 21812.57

# Row 14
$ws.Cells.Item(14, 3).Value = 66
$ws.Cells.Item(14, 4).Value = 14205.35

# Row 15
$ws.Cells.Item(15, 3).Value = 66
$ws.Cells.Item(15, 4).Value = 9713.27

# Row 16
$ws.Cells.Item(16, 3).Value = 66
$ws.Cells.Item(16, 4).Value = 8694.969999999999

# Row 17
$ws.Cells.Item(17, 3).Value = 66
$ws.Cells.Item(17, 4).Value = 7990.55

# Row 18
$ws.Cells.Item(18, 3).Value = 66
$ws.Cells.Item(18, 4).Value = 7445.48

# Row 19
$ws.Cells.Item(19, 3).Value = 66
$ws.Cells.Item(19, 4).Value = 7373.32

# Row 20
$ws.Cells.Item(20, 3).Value = 66
$ws.Cells.Item(20, 4).Value = 7213.18

# Row 21
$ws.Cells.Item(21, 3).Value = 66
$ws.Cells.Item(21, 4).Value = 7109.43

# Row 22
$ws.Cells.Item(22, 3).Value = 66
$ws.Cells.Item(22, 4).Value = 6986.94

# Row 23
$ws.Cells.Item(23, 3).Value = 66
$ws.Cells.Item(23, 4).Value = 6610.19

# Row 24
$ws.Cells.Item(24, 3).Value = 66
$ws.Cells.Item(24, 4).Value = 6552.64

# Row 43
$ws.Cells.Item(43, 3).Value = 10
$ws.Cells.Item(43, 4).Value = 10.74

# Row 45
$ws.Cells.Item(45, 1).Value = 'SMB CI (SMBC)'
$ws.Cells.Item(45, 2).Value = 11
$ws.Cells.Item(45, 3).Value = 12
$ws.Cells.Item(45, 4).Value = 9.83
$ws.Cells.Item(45, 5).Value = -1.64

# Row 46
$ws.Cells.Item(46, 1).Value = 'ONATEL BF (ONTBF)'
$ws.Cells.Item(46, 2).Value = 5
$ws.Cells.Item(46, 3).Value = 8
$ws.Cells.Item(46, 4).Value = 8.5
$ws.Cells.Item(46, 5).Value = -1.28

# Row 47
$ws.Cells.Item(47, 1).Value = 'VIVO ENERGY CI (SHEC)'
$ws.Cells.Item(47, 2).Value = 6
$ws.Cells.Item(47, 3).Value = 4
$ws.Cells.Item(47, 4).Value = 8.18
$ws.Cells.Item(47, 5).Value = 2.63

# Row 48
$ws.Cells.Item(48, 1).Value = 'SETAO CI (STAC)'
$ws.Cells.Item(48, 2).Value = 15
$ws.Cells.Item(48, 3).Value = 11
$ws.Cells.Item(48, 4).Value = 5.93
$ws.Cells.Item(48, 5).Value = 1.83

# Row 49
$ws.Cells.Item(49, 1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws.Cells.Item(49, 4).Value = 4.77
$ws.Cells.Item(49, 5).Value = 3.94
$ws.Cells.Item(49, 6).Value = '🟡 Observer'

# Row 50
$ws.Cells.Item(50, 1).Value = 'SODE CI (SDCC)'
$ws.Cells.Item(50, 3).Value = 8
$ws.Cells.Item(50, 4).Value = 4.02
$ws.Cells.Item(50, 5).Value = -1.83

# Row 51
$ws.Cells.Item(51, 1).Value = 'BICI CI (BICC)'
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = 1
$ws.Cells.Item(51, 4).Value = 3.7
$ws.Cells.Item(51, 5).Value = -1.22

# Row 52
$ws.Cells.Item(52, 1).Value = 'CFAO MOTORS CI (CFAC)'
$ws.Cells.Item(52, 2).Value = 7
$ws.Cells.Item(52, 3).Value = 10
$ws.Cells.Item(52, 4).Value = 3.36
$ws.Cells.Item(52, 5).Value = -6.56

# Row 54
$ws.Cells.Item(54, 3).Value = 65

# Row 55
$ws.Cells.Item(55, 1).Value = 'SOGB CI (SOGC)'
$ws.Cells.Item(55, 2).Value = 7
$ws.Cells.Item(55, 3).Value = 5
$ws.Cells.Item(55, 4).Value = -0.98
$ws.Cells.Item(55, 5).Value = 2.78

# Row 56
$ws.Cells.Item(56, 1).Value = 'TOTALENERGIES MARKETING SN (TTLS)'
$ws.Cells.Item(56, 2).Value = 11
$ws.Cells.Item(56, 3).Value = 14
$ws.Cells.Item(56, 4).Value = -1.75
$ws.Cells.Item(56, 5).Value = 2.34

# Row 57
$ws.Cells.Item(57, 1).Value = 'SONATEL SN (SNTS)'
$ws.Cells.Item(57, 2).Value = 2
$ws.Cells.Item(57, 3).Value = 3
$ws.Cells.Item(57, 4).Value = -3.43
$ws.Cells.Item(57, 5).Value = 0.8

# Row 58
$ws.Cells.Item(58, 1).Value = 'SOLIBRA CI (SLBC)'
$ws.Cells.Item(58, 2).Value = 12
$ws.Cells.Item(58, 3).Value = 13
$ws.Cells.Item(58, 4).Value = -4.41
$ws.Cells.Item(58, 5).Value = 3.91
$ws.Cells.Item(58, 6).Value = '🟡 Observer'

# Row 59
$ws.Cells.Item(59, 1).Value = 'NEI-CEDA CI (NEIC)'
$ws.Cells.Item(59, 3).Value = 8
$ws.Cells.Item(59, 4).Value = -5.59
$ws.Cells.Item(59, 5).Value = 5.88

# Row 60
$ws.Cells.Item(60, 1).Value = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$ws.Cells.Item(60, 3).Value = 9
$ws.Cells.Item(60, 4).Value = -6.92
$ws.Cells.Item(60, 5).Value = -1.35

# Row 64
$ws.Cells.Item(64, 1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws.Cells.Item(64, 4).Value = -23.22
$ws.Cells.Item(64, 5).Value = -2.7

# Row 65
$ws.Cells.Item(65, 1).Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws.Cells.Item(65, 3).Value = 10
$ws.Cells.Item(65, 4).Value = -23.42
$ws.Cells.Item(65, 5).Value = -2.78
